$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.570.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5400"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4009"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.30%  "

$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07792"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.02%  "

$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.86%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.64%  "

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.76%  "

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.640"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.64%  "

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.831.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.083"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.24%  "

$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.576.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.19%  "

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.94%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.455"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.038.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.77%  "

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.65%  "

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.711"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07532"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.84%  "

$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.649"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2256"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  +3.04%  "

$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.987"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.04%  "

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.23%  "

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6315"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.35%  "

$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.407"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.57%  "

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5909"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.39%  "

$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.35%  "

$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.198"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06916"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
